# Completed easy part of merging changes from SWC
#
# For every source file whose SWC-side change could be merged in without any
# manual work, mark its Status (column D) as "merged". A handful of files
# turned out not to need merging at all and are marked "unneeded" instead.
# Files that still require a manual look (e.g. changed on both SH and SWC
# sides, or still under review) are left blank in Status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that are now fully merged (easy, automatic SWC merges).
$mergedRows = @(6,7,14,15,21,25,28,31,33,36,37,39,48,50,57,74,80,81,82,83,84,89,99,106,107,115,116,117,118,126,127,130,131,133,155,156,162,164,165)
foreach ($r in $mergedRows) {
    $ws.Range("D$r").Value = "merged"
}

# Rows that turned out to be unneeded once looked at.
$unneededRows = @(90,91,93,95)
foreach ($r in $unneededRows) {
    $ws.Range("D$r").Value = "unneeded"
}

# Cosmetic touch-ups left behind by the editing session: the separator
# column was nudged a bit wider, and the cursor was left on D4.
$ws.Columns.Item(5).ColumnWidth = 4.45
$ws.Range("D4").Select()

$wb.Application.Calculate()
